$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-6 down to 4-7
$ws.Rows("3:3").Insert()

# Fill in the new row 3 with the new cherry record (Brooks variety)
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 44532
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100103
$ws.Cells.Item(3, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(3, 9).Value = 100103001
$ws.Cells.Item(3, 10).Value = "Cereza"
$ws.Cells.Item(3, 11).Value = "Brooks"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 400
$ws.Cells.Item(3, 14).Value = 27000
$ws.Cells.Item(3, 15).Value = 28000
$ws.Cells.Item(3, 16).Value = 27500
$ws.Cells.Item(3, 17).Value = "$/bandeja 12 kilos"
$ws.Cells.Item(3, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(3, 19).Value = 2292
$ws.Cells.Item(3, 20).Value = 12
